$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.067.26"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "1.901.46"
$ws.Range("E3").Value = "  +1.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.57"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5062"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3926"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09366"
$ws.Range("E9").Value = "  -2.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.138"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.75"
$ws.Range("E11").Value = "  +1.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.367"
$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.78"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("D14").Value = "1.889.15"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.308"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.51"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.83"
$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.216"
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("D23").Value = "28.107.29"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +1.65%  "

$ws.Range("E26").Value = "  +3.52%  "

$ws.Range("D27").Value = "2.106.44"
$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.19"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.27"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.090"
$ws.Range("E31").Value = "  +2.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1070"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.666"
$ws.Range("E35").Value = "  +2.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06651"
$ws.Range("E36").Value = "  -1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02419"
$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.263"
$ws.Range("E40").Value = "  +7.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6387"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.003"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.42"
$ws.Range("E43").Value = "  -0.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5999"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.714"
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.272"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.024"
$ws.Range("E49").Value = "  +1.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.49"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("E51").Value = "  -1.41%  "

# Row 38/39: Algorand and ARBITRUM swap positions with updated values
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.227"
$ws.Range("E38").Value = "  -1.85%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2172"
$ws.Range("E39").Value = "  -0.13%  "
